# ajuste a documento existente
# Cambia el codigo de moneda de "US$" a "USD" en la columna E (Cod. Moneda)
# para todas las filas de datos de la hoja.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2:E19").Value = "USD"
